# Daily attendance processing - 2025-10-02 20:27:10
# Applies updates to "Recorded By" email orderings, PATHOLOGY LAB/MUSEUM
# session renaming/rescheduling, and related Date/Time/Session fields.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / date / time cell updates -------------------------------
$ws.Range('G2').Value = 'servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range('G3').Value = 'shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg '
$ws.Range('G9').Value = 'norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range('C23').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E23').Value = '23/10/2025'
$ws.Range('F23').Value = '12:00:00'
$ws.Range('C24').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E24').Value = '22/10/2025'
$ws.Range('F24').Value = '14:00:00'
$ws.Range('C25').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E25').Value = '04/11/2025'
$ws.Range('F25').Value = '08:00:00'
$ws.Range('C26').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('G29').Value = 'ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Range('G33').Value = 'servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range('G34').Value = 'shaimaa.ahmed@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg '
$ws.Range('G40').Value = 'norhan.mohamed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Range('C54').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('C55').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E55').Value = '23/10/2025'
$ws.Range('F55').Value = '14:00:00'
$ws.Range('C56').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E56').Value = '04/11/2025'
$ws.Range('F56').Value = '10:00:00'
$ws.Range('C57').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('G60').Value = 'ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Range('G64').Value = 'Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range('G65').Value = 'majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range('G71').Value = 'Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg'
$ws.Range('C85').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('C86').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E86').Value = '22/10/2025'
$ws.Range('F86').Value = '10:00:00'
$ws.Range('C87').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E87').Value = '03/11/2025'
$ws.Range('F87').Value = '12:00:00'
$ws.Range('C88').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('G91').Value = 'ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg'
$ws.Range('G95').Value = 'Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range('G96').Value = 'alshimaa.atef@med.asu.edu.egm, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg '
$ws.Range('G102').Value = 'norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Range('C116').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('C117').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E117').Value = '21/10/2025'
$ws.Range('F117').Value = '10:00:00'
$ws.Range('C118').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E118').Value = '02/11/2025'
$ws.Range('F118').Value = '08:00:00'
$ws.Range('C119').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('G122').Value = 'neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range('G126').Value = 'Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range('G127').Value = 'alshimaa.atef@med.asu.edu.egm, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg '
$ws.Range('G133').Value = 'norhan.mohamed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Range('C147').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('F147').Value = '08:00:00'
$ws.Range('C148').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E148').Value = '21/10/2025'
$ws.Range('C149').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E149').Value = '02/11/2025'
$ws.Range('C150').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('F150').Value = '10:00:00'
$ws.Range('G153').Value = 'neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg'
$ws.Range('G157').Value = 'eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Range('G158').Value = 'majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg'
$ws.Range('G164').Value = 'Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg'
$ws.Range('C178').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E178').Value = '22/10/2025'
$ws.Range('F178').Value = '12:00:00'
$ws.Range('C179').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E179').Value = '23/10/2025'
$ws.Range('F179').Value = '14:00:00'
$ws.Range('C180').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('E180').Value = '05/11/2025'
$ws.Range('F180').Value = '08:00:00'
$ws.Range('C181').Value = 'PATHOLOGY LAB/MUSEUM'
$ws.Range('G184').Value = 'ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg'

# --- Numeric-looking text cells (Session numbers) -------------------------
# These must stay text (not be auto-converted to numbers), so we stage the
# text in an unused scratch cell via a text formula, copy it, and use
# PasteSpecial (values only) onto the target cells - this preserves both the
# text data type and the existing cell style/formatting of the destination.
$scratch = $ws.Range("Z1000")
$scratch.Formula = '="3"'
$scratch.Copy() | Out-Null
$ws.Range('D25').PasteSpecial(-4163) | Out-Null
$ws.Range('D56').PasteSpecial(-4163) | Out-Null
$ws.Range('D87').PasteSpecial(-4163) | Out-Null
$ws.Range('D118').PasteSpecial(-4163) | Out-Null
$ws.Range('D149').PasteSpecial(-4163) | Out-Null
$ws.Range('D180').PasteSpecial(-4163) | Out-Null
$scratch.ClearContents() | Out-Null

$scratch.Formula = '="4"'
$scratch.Copy() | Out-Null
$ws.Range('D26').PasteSpecial(-4163) | Out-Null
$ws.Range('D57').PasteSpecial(-4163) | Out-Null
$ws.Range('D88').PasteSpecial(-4163) | Out-Null
$ws.Range('D119').PasteSpecial(-4163) | Out-Null
$ws.Range('D150').PasteSpecial(-4163) | Out-Null
$ws.Range('D181').PasteSpecial(-4163) | Out-Null
$scratch.ClearContents() | Out-Null

$excel.CutCopyMode = 0
